$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# CU 19 (row 33) - task completed ("Hecho"): total estimated hours = 1,
# consumed on Day 9 (column AF) = 1
$ws.Range("F33").Value = "Hecho"
$ws.Range("G33").Value = 1
$ws.Range("AF33").Value = 1

# CU 18 (row 34) - task completed ("Hecho"): total estimated hours = 1,
# consumed on Day 9 (column AF) = 1
$ws.Range("F34").Value = "Hecho"
$ws.Range("G34").Value = 1
$ws.Range("AF34").Value = 1

# Reflect the last selected cell recorded in the saved view state
$ws.Range("AF36").Select() | Out-Null
